# Apply the "gh-pages output" update to 苏州-漫展信息.xlsx
#
# Changes:
#  1. On sheets "展览" (index 1) and "全部类型" (index 4): bump a handful of
#     "想去人数" (F) / "最低票价" (G) counters to newer scraped values.
#  2. Insert a brand-new event row ("苏州·绘时国乙2.0光夜同人only", 2024-11-02)
#     right before the existing "张家港·META萌圆饿了" row on both sheets,
#     pushing the remaining rows down by one. Column A is a running
#     0-based index (row number - 1), so every pushed-down row needs its
#     index bumped by one to match its new row number.

$wb = $excel.ActiveWorkbook

function Update-Counts($ws, $rowMap) {
    foreach ($row in $rowMap.Keys) {
        $vals = $rowMap[$row]
        $ws.Cells.Item($row, 6).Value = $vals[0]
        if ($vals.Count -gt 1) {
            $ws.Cells.Item($row, 7).Value = $vals[1]
        }
    }
}

function Insert-NewEvent($ws, $insertAtRow, $lastRow) {
    # Push everything from $insertAtRow down by one row.
    $ws.Rows.Item($insertAtRow).Insert()

    # Column A: sequential index, bold/centered/bordered like the rest of
    # the column (style gets reset to default by the row insert).
    $a = $ws.Cells.Item($insertAtRow, 1)
    $a.Value = $insertAtRow - 1
    $a.Font.Bold = $true
    $a.HorizontalAlignment = -4108   # xlCenter
    $a.VerticalAlignment = -4160     # xlTop
    $a.Borders.LineStyle = 1         # xlContinuous

    # Column B ("2024-11-02") looks like a date, so Excel would normally
    # auto-convert it to a date serial. Force text storage via a
    # temporary "@" number format, then restore the default (General)
    # formatting via a PasteSpecial(Formats) from a plain cell so the
    # cell's style matches the rest of the column (no explicit style).
    $b = $ws.Cells.Item($insertAtRow, 2)
    $b.NumberFormat = "@"
    $b.Value = "2024-11-02"
    $ws.Cells.Item(1, 4).Copy()
    $b.PasteSpecial(-4122)           # xlPasteFormats

    $ws.Cells.Item($insertAtRow, 3).Value = "苏州·绘时国乙2.0光夜同人only"
    $ws.Cells.Item($insertAtRow, 4).Value = "东亭街588号 南舍别院"
    $ws.Cells.Item($insertAtRow, 5).Value = "2024.11.02 10:30-11.02 20:30"
    $ws.Cells.Item($insertAtRow, 6).Value = 0
    $ws.Cells.Item($insertAtRow, 7).Value = 178
    $ws.Cells.Item($insertAtRow, 8).Value = "https://show.bilibili.com/platform/detail.html?id=91324"
    $ws.Cells.Item($insertAtRow, 9).Value = "//i1.hdslb.com/bfs/openplatform/202408/YauAhbAd1724662566605.jpeg"

    # Every row that got pushed down by the insert keeps its own old data,
    # but the running index in column A (= row number - 1) must advance
    # by one to stay in sync with its new row number.
    for ($r = $insertAtRow + 1; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }
}

# ---------------------------------------------------------------------
# Sheet "展览" (1st sheet)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

Update-Counts $ws1 @{
    3  = @(12723, 54)
    6  = @(52)
    7  = @(34)
    9  = @(4)
    10 = @(12621)
    11 = @(257)
    12 = @(8)
    13 = @(4938)
    14 = @(5880)
    15 = @(172)
    16 = @(82)
    21 = @(13)
}

Insert-NewEvent $ws1 24 26

# ---------------------------------------------------------------------
# Sheet "全部类型" (4th sheet) - same data, shifted down by one row
# because it carries an extra row from the "演出" sheet.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

Update-Counts $ws4 @{
    4  = @(12723, 54)
    7  = @(52)
    8  = @(34)
    10 = @(4)
    11 = @(12621)
    12 = @(257)
    13 = @(8)
    14 = @(4938)
    15 = @(5880)
    16 = @(172)
    17 = @(82)
    22 = @(13)
}

Insert-NewEvent $ws4 25 27
